$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM_rats")

$ws.Cells.Item(8, 1).Value = "Dosso"
$ws.Cells.Item(8, 2).Value = 1990
$ws.Cells.Item(8, 3).Value = "Ocular complications in the old and glucose-intolerant genetically obese (fa/fa) rat"

$ws.Range("A8:B8").WrapText = $true
$ws.Range("A8:B8").VerticalAlignment = -4108
$ws.Rows.Item(8).RowHeight = 17

$ws.Activate()
$ws.Range("C8").Select()
